# Append the latest Nalco PDF run-log entry (2025-08-14 run, SKIPPED - no change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is row 10 (header is row 1), so the new entry goes to row 11.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-08-14 08:41:13 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-14 14:11:13 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""

# Carry over the same direct formatting (centered alignment, etc.) used by the
# other data rows by copying the formats from the previous last row.
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 8))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 8))
$srcRange.Copy()
$dstRange.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
